$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data points for Monssaf (column E) in Sprint 1 and Sprint 2 sections
$ws.Range("E5").Value = 31
$ws.Range("E15").Value = 3

# Re-enter the total formulas as a single fill across B12:G12 so they become
# a shared formula group (as Excel does when you drag-fill / select+enter)
$ws.Range("B12:G12").Formula = "=SUM(B2:B11)"

# Update the active selection to reflect where the user ended up working
$ws.Range("J10").Select() | Out-Null
